$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Pay frequency(Weekly/Semi-monthly)" column (column S) entirely.
# Deleting the whole column shifts Employee Type / Employment status / Leave
# allowance / Works days per year left by one (S<-T, T<-U, U<-V, V<-W).
$ws.Columns("S").Delete()

# Excel's column delete does not always re-anchor the AutoFilter defined name
# onto the now-shifted column, so fix it up explicitly: it should still point
# at the "Employee Type" header cell, which is now S1 instead of T1.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Employees!_FilterDatabase") {
        $n.RefersTo = "=Employees!`$S`$1:`$S`$1"
    }
}
